# Apply updated 'F' column (want-to-go counts) values across sheets
# as described by the commit 'Update gh-pages to output generated at 456a3b4'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1116
$ws.Range("F3").Value = 241
$ws.Range("F4").Value = 228
$ws.Range("F5").Value = 1768
$ws.Range("F6").Value = 659
$ws.Range("F7").Value = 315
$ws.Range("F8").Value = 467
$ws.Range("F9").Value = 4400
$ws.Range("F10").Value = 50
$ws.Range("F11").Value = 450
$ws.Range("F13").Value = 976
$ws.Range("F16").Value = 1892
$ws.Range("F17").Value = 2959
$ws.Range("F18").Value = 1793
$ws.Range("F19").Value = 109
$ws.Range("F21").Value = 165
$ws.Range("F22").Value = 10
$ws.Range("F23").Value = 657
$ws.Range("F24").Value = 918
$ws.Range("F25").Value = 294
$ws.Range("F26").Value = 26
$ws.Range("F27").Value = 2291
$ws.Range("F28").Value = 993
$ws.Range("F29").Value = 2362
$ws.Range("F30").Value = 242
$ws.Range("F31").Value = 699
$ws.Range("F32").Value = 542
$ws.Range("F34").Value = 879
$ws.Range("F35").Value = 407
$ws.Range("F36").Value = 1090
$ws.Range("F37").Value = 890
$ws.Range("F38").Value = 1160
$ws.Range("F39").Value = 10
$ws.Range("F40").Value = 573
$ws.Range("F41").Value = 512
$ws.Range("F43").Value = 275
$ws.Range("F44").Value = 3463

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 12
$ws.Range("F10").Value = 876
$ws.Range("F23").Value = 26

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1116
$ws.Range("F3").Value = 241
$ws.Range("F4").Value = 228
$ws.Range("F5").Value = 12
$ws.Range("F6").Value = 1769
$ws.Range("F7").Value = 659
$ws.Range("F8").Value = 315
$ws.Range("F9").Value = 467
$ws.Range("F10").Value = 4400
$ws.Range("F11").Value = 50
$ws.Range("F17").Value = 2959
$ws.Range("F19").Value = 1793
$ws.Range("F20").Value = 109
$ws.Range("F22").Value = 165
$ws.Range("F23").Value = 876
$ws.Range("F27").Value = 918
$ws.Range("F28").Value = 294
$ws.Range("F29").Value = 2291
$ws.Range("F32").Value = 993
$ws.Range("F33").Value = 2362
$ws.Range("F34").Value = 699
$ws.Range("F35").Value = 542
$ws.Range("F36").Value = 879
$ws.Range("F37").Value = 1090
$ws.Range("F38").Value = 890
$ws.Range("F39").Value = 1160
$ws.Range("F40").Value = 573
$ws.Range("F41").Value = 512
$ws.Range("F45").Value = 26
$ws.Range("F47").Value = 275
$ws.Range("F48").Value = 3463
